$d = $word.ActiveDocument

# --- helper: find the 1-based paragraph index whose Range.Text equals $oldText ---
# (Paragraph.Range.Text includes the trailing paragraph-mark character, so trim it.)
function Find-ParagraphIndex($oldText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        $t = $d.Paragraphs($i).Range.Text
        $t = $t.TrimEnd([char]13, [char]7)
        if ($t -eq $oldText) {
            return $i
        }
    }
    return -1
}

# --- helper: replace a paragraph's whole content via InsertXML, which lets us
#     keep an explicit empty <w:r/> run next to the text run (the COM engine
#     otherwise silently merges/drops an empty run that shares formatting with
#     its neighbour whenever the neighbour's text is edited in place). ---
function Replace-Paragraph($oldText, $innerXml) {
    $idx = Find-ParagraphIndex $oldText
    if ($idx -lt 0) {
        throw "Paragraph not found: $oldText"
    }
    $rng = $d.Paragraphs($idx).Range
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $innerXml + '</w:p>'
    $rng.InsertXML($xml)
}

# 1. Top heading (Heading1) and the bottom bold "title" line both have this exact
#    text; ReplaceAll fixes both in one shot. (The bold run's differing rPr keeps
#    it from merging with its neighbouring empty run, so its <w:r/> survives.)
$d.Content.Find.Execute("Play Monopoly Electric Wins for Free - Game Review", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Play Monopoly Electric Wins for Free", 2)

# 2. "What we like" bullet list items.
Replace-Paragraph "Unique combination of classic board game and online slot" `
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Unique blend of Monopoly board game and online slot</w:t></w:r>'

Replace-Paragraph "Special symbols and bonuses based on Monopoly theme" `
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Bonuses, multipliers, and free spins based on Monopoly theme</w:t></w:r>'

Replace-Paragraph "Neon graphics and interactive board design" `
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Neon graphics replicate classic board game look and feel</w:t></w:r>'

Replace-Paragraph "Medium/low volatility with frequent and consistent wins" `
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Wide range of betting options and automatic turn settings</w:t></w:r>'

# 3. "What we don't like" bullet list items.
Replace-Paragraph "Limited betting options for high rollers" `
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Limited number of symbols</w:t></w:r>'

Replace-Paragraph "Theme may not appeal to players who are not familiar with Monopoly" `
    '<w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr><w:r/><w:r><w:t>Medium/low volatility may not appeal to high-risk players</w:t></w:r>'

# 4. Bottom italic summary line (the bold title line right above it was already
#    fixed by the ReplaceAll in step 1).
Replace-Paragraph "Explore the unique Monopoly Electric Wins video slot game that perfectly blends the classic board game with an online slot. Play for free and experience the interactive board design and special Monopoly-themed symbols and bonuses." `
    '<w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Monopoly Electric Wins and play this unique slot game for free.</w:t></w:r>'
